$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in column H, matching the header style used by G1
$ws.Cells.Item(1, 8).Value = "Label"
$ws.Cells.Item(1, 7).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)

# Fill in Label values: Control rows = 0, MDD rows = 1 (two repeated blocks of rows 2-11 and 12-21)
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(7, 8).Value = 1
$ws.Cells.Item(8, 8).Value = 1
$ws.Cells.Item(9, 8).Value = 1
$ws.Cells.Item(10, 8).Value = 1
$ws.Cells.Item(11, 8).Value = 1

$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(17, 8).Value = 1
$ws.Cells.Item(18, 8).Value = 1
$ws.Cells.Item(19, 8).Value = 1
$ws.Cells.Item(20, 8).Value = 1
$ws.Cells.Item(21, 8).Value = 1

# Update refit values from the refitted NCDE results
$newD3 = [double]"1.26020941904411E-18"
$ws.Cells.Item(3, 4).Value = $newD3
$ws.Cells.Item(3, 5).Value = $newD3

$ws.Cells.Item(11, 6).Value = 0.6558917164802551
